# Updated cryptos list with refreshed prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.622.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.596.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.819.82'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.612.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.604.12'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.54%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '208.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  -3.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0505'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.279.20'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.24'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +14.69%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.599'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.823'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.16'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.732.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0512'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
